$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "2" = @{ "C" = 0.02032577724713747; "D" = 0.02123905908371881; "E" = 0.4213259698304483; "F" = 0.4876327782781402; "G" = 0.002392444783948849; "I" = 0.3601778416014767; "K" = 1.789648023738039; "O" = 1.581338832510752 }
    "3" = @{ "C" = 0.01775018107716164; "D" = 0.01916381650632815; "E" = 0.367491406866236; "F" = 0.4891684632852105; "G" = 0.002395574174240853; "I" = 0.3621903871668408; "K" = 1.565734647479871; "O" = 1.601201355036835 }
    "4" = @{ "C" = 0.01616243825415609; "D" = 0.01788506107078547; "E" = 0.3345335130287737; "F" = 0.4906889377797654; "G" = 0.002397593935047386; "I" = 0.363861789042911; "K" = 1.42777541606057; "O" = 1.615355814903197 }
    "5" = @{ "C" = 0.01551386812438693; "D" = 0.01736285353879907; "E" = 0.3211252723423286; "F" = 0.4914530216281534; "G" = 0.002398441802841388; "I" = 0.364651873610299; "K" = 1.371439425201459; "O" = 1.621613758707895 }
    "6" = @{ "C" = 0.01540608088974693; "D" = 0.01727607585435464; "E" = 0.3189001382564953; "F" = 0.4915886012312498; "G" = 0.002398584090769971; "I" = 0.3647896312779508; "K" = 1.362077933581986; "O" = 1.622682395351859 }
    "7" = @{ "C" = 0.01615369764623864; "D" = 0.01787802281591411; "E" = 0.3343525967992065; "F" = 0.4906986584976281; "G" = 0.002397605269209215; "I" = 0.3638720039411254; "K" = 1.427016116684854; "O" = 1.6154382317402 }
    "8" = @{ "C" = 0.01943904345171177; "D" = 0.02052448209887103; "E" = 0.4027425350159177; "F" = 0.4880419556984066; "G" = 0.002393503443098704; "I" = 0.3607809879871127; "K" = 1.712543151808347; "O" = 1.587779277709274 }
    "9" = @{ "C" = 0.02583021528175777; "D" = 0.02567665227438454; "E" = 0.5377174667796822; "F" = 0.4874494698667888; "G" = 0.002386236102483437; "I" = 0.3582032168834957; "K" = 2.268581320764724; "O" = 1.549203012412619 }
    "10" = @{ "C" = 0.03049327248243117; "D" = 0.02943755942490611; "E" = 0.6375500115541968; "F" = 0.4898778973921623; "G" = 0.002381364900722422; "I" = 0.3584705881495367; "K" = 2.674637379837236; "O" = 1.530577271541091 }
    "11" = @{ "C" = 0.03260733512811953; "D" = 0.03114291724421037; "E" = 0.683142299547967; "F" = 0.4916150734416576; "G" = 0.002379249408727952; "I" = 0.3590696337742472; "K" = 2.858809868285618; "O" = 1.524249945455779 }
    "12" = @{ "C" = 0.03340681471956941; "D" = 0.03178787105852621; "E" = 0.7004348704457755; "F" = 0.4923646880601567; "G" = 0.002378462686832079; "I" = 0.3593657785753948; "K" = 2.928470612674744; "O" = 1.522165437309695 }
    "13" = @{ "C" = 0.03323468065406132; "D" = 0.0316490061365613; "E" = 0.6967093363857941; "F" = 0.4921991497209248; "G" = 0.002378631483847204; "I" = 0.359298906072425; "K" = 2.913471592874373; "O" = 1.522600471347346 }
    "14" = @{ "C" = 0.03267313044707976; "D" = 0.03119599479208546; "E" = 0.6845644008646588; "F" = 0.4916749004168679; "G" = 0.002379184397188539; "I" = 0.359092605214208; "K" = 2.86454254809064; "O" = 1.524072191462494 }
    "15" = @{ "C" = 0.03232902430035267; "D" = 0.03091840331019569; "E" = 0.677128957343399; "F" = 0.4913657603389083; "G" = 0.002379524941759267; "I" = 0.3589752840150169; "K" = 2.834561406717455; "O" = 1.525014319145725 }
    "16" = @{ "C" = 0.03035496576671903; "D" = 0.02932599641400202; "E" = 0.6345742142698469; "F" = 0.4897771671595947; "G" = 0.002381505167782005; "I" = 0.3584411012620023; "K" = 2.662590053372185; "O" = 1.531034206698422 }
    "17" = @{ "C" = 0.02914207804269608; "D" = 0.02834767118086035; "E" = 0.6085153738775659; "F" = 0.4889652017213493; "G" = 0.002382745645154062; "I" = 0.3582361355644181; "K" = 2.556949533697377; "O" = 1.535278871788904 }
    "18" = @{ "C" = 0.02844378293032435; "D" = 0.027784449096103; "E" = 0.5935436374164169; "F" = 0.4885576619796481; "G" = 0.002383468593226979; "I" = 0.3581631427284577; "K" = 2.496136835243931; "O" = 1.537922068243546 }
    "19" = @{ "C" = 0.02820723754516052; "D" = 0.02759366443708444; "E" = 0.588477253406225; "F" = 0.4884298665677989; "G" = 0.002383714998026524; "I" = 0.35814612017775; "K" = 2.47553802506161; "O" = 1.538851574019787 }
    "20" = @{ "C" = 0.02927126205291586; "D" = 0.02845186917772224; "E" = 0.6112876462169936; "F" = 0.4890454747088313; "G" = 0.002382612615758424; "I" = 0.3582533030651547; "K" = 2.568200455008025; "O" = 1.534806115272829 }
    "21" = @{ "C" = 0.03283810068771231; "D" = 0.03132907797594697; "E" = 0.6881308919734579; "F" = 0.4918263874992519; "G" = 0.002379021603751088; "I" = 0.3591513148518573; "K" = 2.878916427528338; "O" = 1.523631433465567 }
    "22" = @{ "C" = 0.03516298224225523; "D" = 0.03320465624219082; "E" = 0.7385152336740219; "F" = 0.4941792492117756; "G" = 0.002376758388999811; "I" = 0.3601424476197295; "K" = 3.081511544790715; "O" = 1.518145254920057 }
    "23" = @{ "C" = 0.03392273340070062; "D" = 0.03220408021786625; "E" = 0.711608538389072; "F" = 0.492874219586497; "G" = 0.002377958673434809; "I" = 0.3595762601598054; "K" = 2.973427251444207; "O" = 1.520906079744123 }
    "24" = @{ "C" = 0.02921286104300691; "D" = 0.02840476365790323; "E" = 0.6100342714326672; "F" = 0.4890089987349313; "G" = 0.002382672727873958; "I" = 0.3582454020006978; "K" = 2.563114157945506; "O" = 1.535019216761242 }
    "25" = @{ "C" = 0.02410686264765616; "D" = 0.02428703176831704; "E" = 0.5010963853915484; "F" = 0.4871103883826322; "G" = 0.00238811953378284; "I" = 0.3585237448484833; "K" = 2.118584186418104; "O" = 1.557944780031164 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Range("$c$r").Value = $data[$r][$c]
    }
}
